$d = $word.ActiveDocument

# --- Insertion 1: version/date line + blank line at the very start of the document ---
$startXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/><w:jc w:val="right"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="16"/><w:szCs w:val="16"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">version 1.01 [ August 22nd, 2015 ]</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/><w:jc w:val="right"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p>'
$d.Range(0, 0).InsertXML($startXml)

# --- Insertion 2: card list contents inside the "Card List" table's second row ---
$cardListTable = $null
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)
    if ($t.Rows.Item(1).Cells.Item(1).Range.Text.TrimEnd([char]13, [char]7) -eq "Card List") {
        $cardListTable = $t
        break
    }
}

$targetCell = $cardListTable.Rows.Item(2).Cells.Item(1)
$targetPara = $targetCell.Range.Paragraphs.Item(1)
$insertPoint = $d.Range($targetPara.Range.Start, $targetPara.Range.Start)

$cardListXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">[ heat cards ] ( 24 )</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">12x coal</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">5x coal duo</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">3x coal trio</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">4x fire log</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">[ action cards ] ( 34 )</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">8x roasting stick</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2x fishing stick</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">12x fire stick</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">4x shovel</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2x balloon of propane</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2x spaghetti</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">4x bucket of water</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rtl w:val="0"/></w:rPr></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">[ play immediately cards ] ( 10 )</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2x lightning strike</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2x rainstorm</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">2x wandering coals</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">1x santa coals</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">1x mrs. coals</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">1x hungry uncle</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:contextualSpacing w:val="0"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Play" w:cs="Play" w:eastAsia="Play" w:hAnsi="Play"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">1x guilt trip</w:t></w:r></w:p>'
$insertPoint.InsertXML($cardListXml)

Write-Output "done"
